# Auto-generated: apply per-cell text updates to match the target diff.
# Column D holds numeric-looking price strings stored as text in the source
# data; a leading apostrophe forces Excel to keep them as text (preserving
# exact formatting like trailing zeros, e.g. "0.8380") instead of silently
# coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.749.88'
$ws.Range('E2').Value = '  -2.76%  '
$ws.Range('D3').Value = '''1.774.65'
$ws.Range('E3').Value = '  -3.10%  '
$ws.Range('D4').Value = '''1.007'
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').Value = '''1.006'
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('D6').Value = '''303.04'
$ws.Range('E6').Value = '  -3.08%  '
$ws.Range('D7').Value = '''0.4355'
$ws.Range('E7').Value = '  +1.56%  '
$ws.Range('D8').Value = '''0.3625'
$ws.Range('E8').Value = '  -1.00%  '
$ws.Range('D9').Value = '''0.07187'
$ws.Range('E9').Value = '  -1.18%  '
$ws.Range('D10').Value = '''0.8380'
$ws.Range('E10').Value = '  -3.34%  '
$ws.Range('D11').Value = '''20.22'
$ws.Range('E11').Value = '  -2.02%  '
$ws.Range('D12').Value = '''1.777.42'
$ws.Range('E12').Value = '  -4.57%  '
$ws.Range('D13').Value = '''5.254'
$ws.Range('E13').Value = '  -2.73%  '
$ws.Range('D14').Value = '''6.357'
$ws.Range('E14').Value = '  -2.87%  '
$ws.Range('D15').Value = '''0.06808'
$ws.Range('E15').Value = '  -1.93%  '
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('D17').Value = '''79.30'
$ws.Range('E17').Value = '  -1.65%  '
$ws.Range('D18').Value = '''0.000008686'
$ws.Range('E18').Value = '  -2.54%  '
$ws.Range('E19').Value = '  +0.45%  '
$ws.Range('D20').Value = '''14.97'
$ws.Range('E20').Value = '  -3.02%  '
$ws.Range('D21').Value = '''27.051.01'
$ws.Range('E21').Value = '  -2.35%  '
$ws.Range('D22').Value = '''5.011'
$ws.Range('E22').Value = '  -2.75%  '
$ws.Range('D23').Value = '''11.05'
$ws.Range('E23').Value = '  +1.64%  '
$ws.Range('D24').Value = '''2.056.46'
$ws.Range('E24').Value = '  -1.51%  '
$ws.Range('D25').Value = '''1.908'
$ws.Range('E25').Value = '  -3.71%  '
$ws.Range('D26').Value = '''153.68'
$ws.Range('E26').Value = '  -0.60%  '
$ws.Range('D27').Value = '''18.13'
$ws.Range('E27').Value = '  -3.97%  '
$ws.Range('D28').Value = '''114.83'
$ws.Range('E28').Value = '  +0.63%  '
$ws.Range('D29').Value = '''5.034'
$ws.Range('E29').Value = '  -1.78%  '
$ws.Range('D30').Value = '''1.636'
$ws.Range('E30').Value = '  -10.97%  '
$ws.Range('D31').Value = '''0.08971'
$ws.Range('E31').Value = '  +1.32%  '
$ws.Range('D32').Value = '''0.7203'
$ws.Range('E32').Value = '  -4.30%  '
$ws.Range('D33').Value = '''2.841'
$ws.Range('E33').Value = '  -5.60%  '
$ws.Range('D34').Value = '''4.325'
$ws.Range('E34').Value = '  -4.89%  '
$ws.Range('D35').Value = '''1.091'
$ws.Range('E35').Value = '  -3.80%  '
$ws.Range('D36').Value = '''1.006'
$ws.Range('E36').Value = '  +0.44%  '
$ws.Range('D37').Value = '''1.072'
$ws.Range('E37').Value = '  -1.50%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '''0.05095'
$ws.Range('E38').Value = '  -4.33%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.01887'
$ws.Range('E39').Value = '  -2.64%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '''0.4920'
$ws.Range('E40').Value = '  -3.29%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = '''0.1609'
$ws.Range('E41').Value = '  -3.26%  '
$ws.Range('D42').Value = '''2.582'
$ws.Range('E42').Value = '  -7.56%  '
$ws.Range('D43').Value = '''6.129'
$ws.Range('E43').Value = '  -6.63%  '
$ws.Range('D44').Value = '''7.932'
$ws.Range('E44').Value = '  -4.85%  '
$ws.Range('D45').Value = '''104.74'
$ws.Range('E45').Value = '  -1.20%  '
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').Value = '''10.05'
$ws.Range('E47').Value = '  -3.32%  '
$ws.Range('D48').Value = '''0.06229'
$ws.Range('E48').Value = '  -4.22%  '
$ws.Range('D49').Value = '''0.4488'
$ws.Range('E49').Value = '  -4.28%  '
$ws.Range('D50').Value = '''1.577'
$ws.Range('E50').Value = '  -2.53%  '
$ws.Range('D51').Value = '''1.722'
$ws.Range('E51').Value = '  -0.12%  '
